$p = $ppt.ActivePresentation
$newStyleId = "{E926617A-ABBB-4531-9791-336283F6A914}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newStyleId)
        }
    }
}
